# Generate Report for Handoff
# Adds a new file (d8893c56-fe19-464e-94ef-f52a94eba1ad.md) as a row on each
# of the three sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1) -> table3 "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A5").Value = "d8893c56-fe19-464e-94ef-f52a94eba1ad.md"
$wsOverview.Range("B5").Value = "e2e\d8893c56-fe19-464e-94ef-f52a94eba1ad.md"
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("D5").Value = "'"
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2017-02-21 11:08:44"
$wsOverview.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/2f2cc2edf61b9592c5a5679de992c9f8255c2241/e2e/d8893c56-fe19-464e-94ef-f52a94eba1ad.md", "", "", "e2e\d8893c56-fe19-464e-94ef-f52a94eba1ad.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) -> table1 "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A5").Value = "d8893c56-fe19-464e-94ef-f52a94eba1ad.md"
$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = "e2e"
$wsZh.Range("E5").Value = "ht"
$wsZh.Range("F5").Value = "'False"
$wsZh.Range("G5").Value = "d8893c56-fe19-464e-94ef-f52a94eba1ad.ce58ac674dba2c9d903ff8d34ca36f55a512b5c0.zh-cn.xlf"
$wsZh.Range("H5").Value = "2017-02-21 11:08:27"
$wsZh.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I5").Value = "'"
$wsZh.Range("J5").Value = "'"
$wsZh.Range("K5").Value = "'"
$wsZh.Range("L5").Value = "0001-01-01 00:00:00"
$wsZh.Range("L5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("M5").Value = "'"
$wsZh.Range("N5").Value = "'"
$wsZh.Range("O5").Value = "'True"
$wsZh.Range("P5").Value = "'"
$wsZh.Range("Q5").Value = "'False"
$wsZh.Range("R5").Value = "'"

$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test4-zhcn/blob/6641740da21c6f1bddb1c3fc20eb0b5cff31430c/e2e/d8893c56-fe19-464e-94ef-f52a94eba1ad.md", "", "", "d8893c56-fe19-464e-94ef-f52a94eba1ad.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3) -> table2 "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A5").Value = "d8893c56-fe19-464e-94ef-f52a94eba1ad.md"
$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = "e2e"
$wsDe.Range("E5").Value = "ht"
$wsDe.Range("F5").Value = "'False"
$wsDe.Range("G5").Value = "d8893c56-fe19-464e-94ef-f52a94eba1ad.ce58ac674dba2c9d903ff8d34ca36f55a512b5c0.de-de.xlf"
$wsDe.Range("H5").Value = "2017-02-21 11:08:44"
$wsDe.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I5").Value = "'"
$wsDe.Range("J5").Value = "'"
$wsDe.Range("K5").Value = "'"
$wsDe.Range("L5").Value = "0001-01-01 00:00:00"
$wsDe.Range("L5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("M5").Value = "'"
$wsDe.Range("N5").Value = "'"
$wsDe.Range("O5").Value = "'True"
$wsDe.Range("P5").Value = "'"
$wsDe.Range("Q5").Value = "'False"
$wsDe.Range("R5").Value = "'"

$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test4-dede/blob/ac0e4da03c29d89435978d38f244eaa0523d7eba/e2e/d8893c56-fe19-464e-94ef-f52a94eba1ad.md", "", "", "d8893c56-fe19-464e-94ef-f52a94eba1ad.md") | Out-Null

Write-Host "Generate Report for Handoff: added row 5 to Overview, zh-cn, de-de"
